# Split the single run that contains the long sentence about the
# clustering algorithm into several runs, inserting extra clarifying
# text along the way (per the commit's diff):
#
#   "Após a formação dos grupos, o algoritmo retorna o identificador do
#    grupo que foi atribuído a cada amostra do conjunto auxiliar. Nesse
#    momento, as características d"
#
# becomes 7 runs:
#
#   "Após a formação dos grupos, o algoritmo retorna"
#   " um vetor com"
#   " o identificador do grupo"
#   " (valor armazenado no índice i)"
#   " que foi atribuído a cada amostra do conjunto auxiliar"
#   " (índice i)"
#   ". Nesse momento, as características d"
#
# A plain Find-and-Replace (or Range.Text / Range.InsertAfter without
# tracked changes) would let Word's writer silently re-coalesce all the
# newly written text -- and any untouched neighboring runs it touches --
# back into a single <w:r>. To force Word to keep the new text as
# distinct runs, we perform the edit with revision tracking turned on
# (each tracked insertion becomes its own <w:ins>/<w:r>), then accept
# all the revisions, which leaves the separate runs in place.

$d = $word.ActiveDocument

$oldText = "Após a formação dos grupos, o algoritmo retorna o identificador do grupo que foi atribuído a cada amostra do conjunto auxiliar. Nesse momento, as características d"

$segments = @(
  "Após a formação dos grupos, o algoritmo retorna",
  " um vetor com",
  " o identificador do grupo",
  " (valor armazenado no índice i)",
  " que foi atribuído a cada amostra do conjunto auxiliar",
  " (índice i)",
  ". Nesse momento, as características d"
)

$originalUserName = $word.UserName
$originalTrackRevisions = $d.TrackRevisions

$d.TrackRevisions = $true

# Keep using the very Range object Find.Execute was called on: on a
# successful match it collapses to the found text, giving us Start/End
# for that occurrence. (Re-fetching $d.Content afterwards would just
# give a fresh range over the whole document again.)
$searchRange = $d.Content
$found = $searchRange.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $searchRange.Start
    $end = $searchRange.End

    # Overwrite the whole old run's text with the first new segment.
    $word.UserName = "Segment0"
    $target = $d.Range($start, $end)
    $target.Text = $segments[0]

    # The deleted original text is still physically present in the
    # story (as a tracked <w:del>) right after the first inserted
    # segment, so the insertion point for the remaining segments sits
    # at start + len(segment0) + len(oldText), and then advances by
    # each segment's length as it's inserted.
    $pos = $start + $segments[0].Length + $oldText.Length

    for ($i = 1; $i -lt $segments.Length; $i++) {
        # Give every insertion a distinct "author" so Word treats each
        # one as a separate revision instead of merging adjacent
        # same-author tracked insertions back into one run.
        $word.UserName = "Segment$i"
        $insertionPoint = $d.Range($pos, $pos)
        $insertionPoint.InsertAfter($segments[$i])
        $pos = $pos + $segments[$i].Length
    }

    $d.Revisions.AcceptAll()
}

$d.TrackRevisions = $originalTrackRevisions
$word.UserName = $originalUserName
